$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add two new columns (P, Q) to the header row, copying the header
#     formatting (bold font, thin border, centered) from the existing O1 cell ---
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Cells.Item(1, 16).Value = 14
$ws.Cells.Item(1, 17).Value = 15

# --- Rewrite the data grid (rows 2-25, columns B-Q). Column A (the index
#     column) and columns J-N are unchanged; B-I, O, P and Q get new values
#     and P/Q are brand-new columns. ---
$data = @(
    @(24.49266632314295,18.96828366782543,4.333588960490924,29.21382029618269,18.59124872237463,2.071584505415577,3.159216038339132,3.138252493947185,0,0,0,0,0,0,0,14.5687144452271),
    @(22.88573587628169,17.85725570160992,4.196859081625459,27.32102133038484,18.00587086705919,2.076468730828436,2.90028629522092,2.942540537153831,0,0,0,0,0,0,0,14.3407332244305),
    @(21.83970279736576,17.14498019991685,4.110482788901874,26.09414241005118,17.64924845131027,2.07955589531849,2.735997085624268,2.819518945045894,0,0,0,0,0,0,0,14.2082728333152),
    @(21.39827057864807,16.86313572588381,4.076836119671445,25.57764346556439,17.49084442949468,2.080842831687721,2.667542724177149,2.769161091993968,0,0,0,0,0,0,0,14.14526399565041),
    @(21.32384007689132,16.8353225505169,4.073836145152685,25.49061016046219,17.44754621655228,2.081065407239584,2.655917759783867,2.761577121908168,0,0,0,0,0,0,0,14.12160564649061),
    @(21.83321247646425,17.19321626836599,4.117147111988736,26.08648815999262,17.60102890494004,2.079593170369083,2.734610986851889,2.820713870448712,0,0,0,0,0,0,0,14.17126314780701),
    @(23.95016255225594,18.65577656276812,4.296175493917771,28.57366895627258,18.33033639464374,2.073275125169522,3.070445310790362,3.073227600650782,0,0,0,0,0,0,0,14.44169348795232),
    @(27.63402988064382,21.197725204966,4.614753228655371,32.94140957282742,19.85430217178764,2.06152701904097,3.688970248645137,3.5447492924246,0,0,0,0,0,0,0,15.09616439954949),
    @(30.05482127450283,22.928038186333,4.874742813704069,34.9646762972166,20.76645788651782,2.053507625778628,4.0828560900355,3.872324092016694,0,0,0,0,0,0,0,15.44422833104606),
    @(31.06280360374938,23.67909720321394,5.319755378876647,28.3645749766292,19.49143900051977,2.051789214662098,4.451631542486619,3.959039045160723,0,0,0,0,0,0,0,14.22126316644512),
    @(31.42548808711483,23.9023462235017,5.642827013362487,22.41126104158498,18.26442248682248,2.051750605528799,5.267815255078704,3.968982522168912,0,0,0,0,0,0,0,13.1551535732224),
    @(31.32574147658034,23.8023350085204,5.904668838461662,16.48110235038558,16.90881854901085,2.053024761394272,6.295261935674344,3.922986396371375,0,0,0,0,0,0,0,12.08345473812526),
    @(31.04587492091814,23.59796749859168,6.064911838226821,12.39948689670872,15.87950193874549,2.054486513402981,7.098170924744764,3.867096236766066,0,0,0,0,0,0,0,11.32485027251851),
    @(30.8771560147254,23.49405501822589,6.093011043100423,11.41723285621635,15.5867530759949,2.055159410205087,7.282752320104646,3.841502725540322,0,0,0,0,0,0,0,11.1299896727916),
    @(29.90608954675489,22.82702198544147,5.945605588702913,11.28027543478703,15.37602402775759,2.058301267800048,7.001793318481235,3.713436795145757,0,0,0,0,0,0,0,11.15774339611703),
    @(29.29981657985098,22.41983437938699,5.742303216938881,13.38635553275821,15.7802112669253,2.059961030497345,6.307931783648366,3.644269112860384,0,0,0,0,0,0,0,11.58299135364639),
    @(28.95452436606275,22.16444665626625,5.465780429740172,17.95867208990727,16.78554409291046,2.060395421184307,5.266155076833628,3.619276282305147,0,0,0,0,0,0,0,12.42595252609114),
    @(28.85164769548753,22.13903500602907,5.173506130908312,24.24597306498533,18.11461374072026,2.059659617645809,4.278060908810593,3.643304718169249,0,0,0,0,0,0,0,13.49383433671839),
    @(29.43636915149335,22.61293314624111,4.829660403423384,34.40385939798826,20.3946872154569,2.0556477065682,3.977166495449919,3.792827574933824,0,0,0,0,0,0,0,15.24019517860266),
    @(31.2128873691445,23.87180592760401,4.963924874935277,37.23156527326921,21.41888838623526,2.049227471009077,4.331854164063795,4.05173183932555,0,0,0,0,0,0,0,15.78587279155189),
    @(32.32256361032284,24.6138801197208,5.059435747161129,38.57128553915183,22.03675548966449,2.045186385304464,4.540833089149837,4.212966268023601,0,0,0,0,0,0,0,16.11096486298121),
    @(31.73534804202202,24.17410756943924,5.001185173529103,37.86180514530616,21.75145750447885,2.047321382313501,4.430018046149989,4.124763455075183,0,0,0,0,0,0,0,15.97437081622241),
    @(29.41053052578191,22.52351586716033,4.786763194111646,35.06520157565711,20.60081722598498,2.055542008347234,4.002604939480796,3.79065573902421,0,0,0,0,0,0,0,15.41479316140643),
    @(26.68629855093528,20.61659098146711,4.543036311527937,31.81296454054986,19.36590674883738,2.064667127147425,3.52510738843625,3.422478011629868,0,0,0,0,0,0,0,14.85062811068665)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $rowVals.Length; $j++) {
        $col = $j + 2
        $ws.Cells.Item($row, $col).Value = $rowVals[$j]
    }
}
